$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = "0/3"
$ws.Range("B4").Value = "0/2"
